$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells that are being updated stay as text (matches original inlineStr formatting)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.712.72'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.902.40'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '311.26'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').Value = '0.5169'
$ws.Range('E7').Value = '  +4.09%  '
$ws.Range('D8').Value = '0.3789'
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = '0.07237'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('D10').Value = '21.20'
$ws.Range('E10').Value = '  +2.99%  '
$ws.Range('D11').Value = '0.9003'
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.897.17'
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.07657'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').Value = '5.442'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').Value = '92.16'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').Value = '14.40'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '1.0000'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '27.739.81'
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').Value = '5.140'
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('D22').Value = '2.159.39'
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').Value = '152.86'
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('D26').Value = '1.856'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('D27').Value = '18.27'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').Value = '2.154'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').Value = '113.68'
$ws.Range('E29').Value = '  -1.09%  '
$ws.Range('D30').Value = '4.820'
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.831'
$ws.Range('E31').Value = '  +4.30%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.09082'
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.05283'
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '3.185'
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('D35').Value = '1.226'
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.7747'
$ws.Range('E36').Value = '  +1.54%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.02084'
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '2.599'
$ws.Range('E38').Value = '  +2.49%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '3.073'
$ws.Range('E39').Value = '  +2.66%  '
$ws.Range('D40').Value = '0.5575'
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.092'
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('D42').Value = '6.676'
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('D43').Value = '117.21'
$ws.Range('E43').Value = '  +4.00%  '
$ws.Range('D44').Value = '8.500'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').Value = '0.1516'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').Value = '0.4809'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '10.40'
$ws.Range('E47').Value = '  -1.55%  '
$ws.Range('D48').Value = '0.9993'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').Value = '1.607'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('D50').Value = '66.70'
$ws.Range('E50').Value = '  -0.94%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = '36.92'
$ws.Range('E51').Value = '  +0.18%  '
